$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly added helper/tutor names and time slots ("det and algebra cofactor add")
$ws.Range("H4").Value = "邓涵朵"
$ws.Range("H5").Value = "向东伟"
$ws.Range("K5").Value = "18日午间十分钟"
$ws.Range("H7").Value = "陶昱"
$ws.Range("K7").Value = "19日午间十分钟"

# Move the active selection to H6:J6 as captured by the saved view state
$ws.Range("H6:J6").Select()
